$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tardigrade")

# K3 corresponds to Tardigrade.Framework.AspNet, release 12.1.0 -> set version 5.2.0
$ws.Range("K3").Value = "5.2.0"
$ws.Range("C3").Copy()
$ws.Range("K3").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# K7 corresponds to Tardigrade.Framework.EntityFramework, release 12.1.0 -> set version 11.0.0
$ws.Range("K7").Value = "11.0.0"

# Reflect the final selection being on K7
$ws.Range("K7").Select()
